$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.222.40'
$ws.Range('E2').Value = '  -3.29%  '
$ws.Range('D3').Value = '2.289.02'
$ws.Range('E3').Value = '  -5.21%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '545.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.68'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.19%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.570'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.06%  '
$ws.Range('D9').Value = '2.288.27'
$ws.Range('E9').Value = '  -5.18%  '
$ws.Range('E10').Value = '  -4.55%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.47'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.149'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.331'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.57'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.19%  '
$ws.Range('D15').Value = '2.709.23'
$ws.Range('E15').Value = '  -4.81%  '
$ws.Range('D16').Value = '58.279.06'
$ws.Range('E16').Value = '  -3.09%  '
$ws.Range('E17').Value = '  -4.72%  '
$ws.Range('D18').Value = '2.292.83'
$ws.Range('E18').Value = '  -6.88%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.57'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.26'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.84%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '312.12'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.42'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.46%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.88'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.167'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.21%  '
$ws.Range('E26').Value = '  +0.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.96'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -7.52%  '
$ws.Range('E28').Value = '  -6.79%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.74'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.68%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '170.29'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('D31').Value = '0.0₃0716'
$ws.Range('E31').Value = '  -7.37%  '
$ws.Range('E32').Value = '  -0.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.71'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -7.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.377'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.68'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.14%  '
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('E38').Value = '  -8.25%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.92'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '38.00'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.50%  '
$ws.Range('E41').Value = '  -6.94%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '290.20'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -10.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '139.61'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.18%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.40'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0947'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0499'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.60%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.552'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.21'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -7.98%  '
$ws.Range('E49').Value = '  -4.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '10.99'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '16.47'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.46%  '
